# Commit: Input mean years 2009,11,13
# Updates the quintile output values (columns A-D) on rows whose underlying
# per-year input means were recalculated (years 2009, 2011, 2013).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 199.33698629
$ws.Range("B2").Value = 234.9835616545
$ws.Range("C2").Value = 271.521917789
$ws.Range("D2").Value = 573.3808218950001
$ws.Range("A3").Value = 29.65
$ws.Range("B3").Value = 69.8
$ws.Range("C3").Value = 123.35
$ws.Range("D3").Value = 348.6
$ws.Range("A4").Value = 63.71506849
$ws.Range("B4").Value = 160.547945237
$ws.Range("C4").Value = 269.4465753
$ws.Range("D4").Value = 407.2657533870001
$ws.Range("A5").Value = 57.6
$ws.Range("B5").Value = 145.6
$ws.Range("C5").Value = 257.2
$ws.Range("D5").Value = 469.8
$ws.Range("A6").Value = 103.367123283
$ws.Range("B6").Value = 108.350684929
$ws.Range("C6").Value = 109.6246575275
$ws.Range("D6").Value = 110.898630126
$ws.Range("A7").Value = 208.7
$ws.Range("B7").Value = 316
$ws.Range("C7").Value = 446.6
$ws.Range("D7").Value = 905.2
$ws.Range("A8").Value = 72.27534247124999
$ws.Range("B8").Value = 124.4328767135
$ws.Range("C8").Value = 167.338356161
$ws.Range("D8").Value = 201.983561631
$ws.Range("A9").Value = 58.15000000000001
$ws.Range("B9").Value = 91.2
$ws.Range("C9").Value = 127.15
$ws.Range("D9").Value = 265
$ws.Range("A10").Value = 454.863013717
$ws.Range("B10").Value = 522.630137034
$ws.Range("C10").Value = 608.6630136935
$ws.Range("D10").Value = 812.5753424099998
$ws.Range("A12").Value = 21.59999999975
$ws.Range("B12").Value = 34.3068493175
$ws.Range("C12").Value = 52.17123287875
$ws.Range("D12").Value = 103.517808213
$ws.Range("B13").Value = 15.6
$ws.Range("C13").Value = 30.59999999999999
$ws.Range("D13").Value = 35.8
$ws.Range("A14").Value = 327.62876710125
$ws.Range("B14").Value = 350.4164383965
$ws.Range("C14").Value = 402.0931506715
$ws.Range("D14").Value = 718.1972601889998
$ws.Range("A15").Value = 272.1
$ws.Range("B15").Value = 367.5999999999999
$ws.Range("C15").Value = 485.9
$ws.Range("D15").Value = 981.8000000000003
$ws.Range("A17").Value = 1.2
$ws.Range("C17").Value = 8.4
$ws.Range("D17").Value = 78.8
$ws.Range("A20").Value = 102.9315068535
$ws.Range("B20").Value = 120.334246518
$ws.Range("C20").Value = 133.460273994
$ws.Range("D20").Value = 245.709589081
$ws.Range("A21").Value = 170.8
$ws.Range("B21").Value = 293.6
$ws.Range("C21").Value = 390.6
$ws.Range("D21").Value = 583.5999999999999
$ws.Range("A22").Value = 56.84657534
$ws.Range("B22").Value = 93.41917809
$ws.Range("C22").Value = 123.169863016
$ws.Range("D22").Value = 156.027397262
$ws.Range("A23").Value = 26.6
$ws.Range("B23").Value = 46.59999999999999
$ws.Range("C23").Value = 68.59999999999999
$ws.Range("D23").Value = 219
